$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 22.238109643292994
$ws.Range("C2").Value = -1.9546825443175493
$ws.Range("D2").Value = 0.77861807930133864
$ws.Range("E2").Value = 3.3322129417138484

# Row 3 data values
$ws.Range("B3").Value = 19.973490344849282
$ws.Range("C3").Value = 3.6266888763321674
$ws.Range("D3").Value = -8.8444066458521178
$ws.Range("E3").Value = 15.967788290335697

# Update the selection shown in the sheet view (B1:AY3 -> B1:E3)
$ws.Range("B1:E3").Select()
